$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts the existing table data
# (Team ... Yellow Cards) from A:M to B:N, keeping everything intact.
$ws.Columns("A:A").Insert()

# Populate the new "Season" column.
$ws.Range("A1").Value = "Season"
$ws.Range("A2").Value = "24/25"
$ws.Range("A3").Value = "23/24"

# Add the new 23/24 season row of stats (table data now lives in columns B:N).
$ws.Range("B3").Value = "Newcastle"
$ws.Range("C3").Value = 7
$ws.Range("D3").Value = 60
$ws.Range("E3").Value = 38
$ws.Range("F3").Value = 18
$ws.Range("G3").Value = 14
$ws.Range("H3").Value = 6
$ws.Range("I3").Value = 85
$ws.Range("J3").Value = 62
$ws.Range("K3").Value = 23
$ws.Range("L3").Value = "Alexander Isak"
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = 77

# The inserted column shifted Table1's backing range but left its cached
# column headers stale, so rebuild the table over the new extent B1:N3.
$lo = $ws.ListObjects("Table1")
$lo.Unlist()
$lo2 = $ws.ListObjects.Add(1, $ws.Range("B1:N3"), $null, 1)
$lo2.Name = "Table1"

# Match the author's final selection.
$ws.Range("E12").Select() | Out-Null
